# edit.ps1 - Update QE_holdings model: refresh "as of" date in the
# confidential disclosure note and refresh the Weight/Percent Change
# figures for each holding (rows 2-35) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; remove protection temporarily so the cells
# (which are locked by default) can be edited, and re-protect afterwards.
$ws.Unprotect("")

# Roll the "Model holdings provided as of ..." date forward from
# 2021-07-09 to 2021-07-13 in the confidential disclosure note.
$noteCell = $ws.Range("A38")
$noteCell.Value = $noteCell.Value2 -replace "2021-07-09", "2021-07-13"

# Refresh the Weight (column D) and Percent Change (column E) values
# for every holding row (2-35, including the "Total" row 35).
$ws.Range("D2").Value = [double]"0.09402055042829777"
$ws.Range("E2").Value = [double]"0.007889273356401327"
$ws.Range("D3").Value = [double]"0.07845277281412828"
$ws.Range("E3").Value = [double]"0.01319774989182187"
$ws.Range("D4").Value = [double]"0.05259818230707788"
$ws.Range("E4").Value = [double]"-0.01107689825335145"
$ws.Range("D5").Value = [double]"0.05150733457892528"
$ws.Range("E5").Value = [double]"-0.003213023455071462"
$ws.Range("D6").Value = [double]"0.04710476254329904"
$ws.Range("E6").Value = [double]"0.01883381679068385"
$ws.Range("D7").Value = [double]"0.04246271887855422"
$ws.Range("E7").Value = [double]"-0.01487341772151907"
$ws.Range("D8").Value = [double]"0.03592088581588183"
$ws.Range("E8").Value = [double]"0.002882445826163149"
$ws.Range("D9").Value = [double]"0.0383561977110586"
$ws.Range("E9").Value = [double]"-0.00123908425772945"
$ws.Range("D10").Value = [double]"0.0336766703004683"
$ws.Range("E10").Value = [double]"0.003784362727597079"
$ws.Range("D11").Value = [double]"0.03530883989786882"
$ws.Range("E11").Value = [double]"0.006009037592539235"
$ws.Range("D12").Value = [double]"0.03496766708662608"
$ws.Range("E12").Value = [double]"-0.003029788197984162"
$ws.Range("D13").Value = [double]"0.0310339954943123"
$ws.Range("E13").Value = [double]"-0.01895151365985726"
$ws.Range("D14").Value = [double]"0.03185665762954772"
$ws.Range("E14").Value = [double]"-0.01457698762976334"
$ws.Range("D15").Value = [double]"0.03187334850505628"
$ws.Range("E15").Value = [double]"0.02170092661625289"
$ws.Range("D16").Value = [double]"0.03129624026713322"
$ws.Range("E16").Value = [double]"-0.003959214665364863"
$ws.Range("D17").Value = [double]"0.02909728814292934"
$ws.Range("E17").Value = [double]"-0.001239609158524058"
$ws.Range("D18").Value = [double]"0.02909728814292935"
$ws.Range("E18").Value = [double]"-0.02502552136502867"
$ws.Range("D19").Value = [double]"0.02389822187011076"
$ws.Range("E19").Value = [double]"-0.007552351527634782"
$ws.Range("D20").Value = [double]"0.02065000775842815"
$ws.Range("E20").Value = [double]"-0.00106856634016006"
$ws.Range("D21").Value = [double]"0.02163095031479981"
$ws.Range("E21").Value = [double]"-0.0047408860552558"
$ws.Range("D22").Value = [double]"0.02157705858964082"
$ws.Range("E22").Value = [double]"0.009911894273127997"
$ws.Range("D23").Value = [double]"0.02086331148950363"
$ws.Range("E23").Value = [double]"0.002467830072272026"
$ws.Range("D24").Value = [double]"0.01853083236123128"
$ws.Range("E24").Value = [double]"-0.007373595505617891"
$ws.Range("D25").Value = [double]"0.0221711405992675"
$ws.Range("E25").Value = [double]"-0.0001148369315573383"
$ws.Range("D26").Value = [double]"0.02024631488807707"
$ws.Range("E26").Value = [double]"-0.002766599597585517"
$ws.Range("D27").Value = [double]"0.01957613380197917"
$ws.Range("E27").Value = [double]"0.00169077587826405"
$ws.Range("D28").Value = [double]"0.01848599331431422"
$ws.Range("E28").Value = [double]"-0.001346687989226436"
$ws.Range("D29").Value = [double]"0.02060021802809753"
$ws.Range("E29").Value = [double]"-0.001421332344580306"
$ws.Range("D30").Value = [double]"0.01160581640234968"
$ws.Range("E30").Value = [double]"-0.01279707495429616"
$ws.Range("D31").Value = [double]"0.008570906021742554"
$ws.Range("E31").Value = [double]"-0.005875169158662708"
$ws.Range("D32").Value = [double]"0.00760014772839306"
$ws.Range("E32").Value = [double]"0.00627198451545663"
$ws.Range("D33").Value = [double]"0.008489714813760242"
$ws.Range("E33").Value = [double]"0.008013995334888335"
$ws.Range("D34").Value = [double]"0.006871831474210263"
$ws.Range("E34").Value = [double]"0.004569593676670536"
$ws.Range("D35").Value = [double]"1"
$ws.Range("E35").Value = [double]"1.753956409378254E-05"

# Restore sheet protection so the workbook ends up protected again,
# matching its state before the edit.
$ws.Protect("")

"Updated note date and " + 34 + " holding rows."
